$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing "KCEPS" row label (A7) which was blank before.
$ws.Range("A7").Value = "KCEPS"

# Fill in the Low / Moderate / High counts for each source row (B:D).
$data = @{
    2  = @(4440, 20468, 5345)
    3  = @(248, 2307, 1471)
    4  = @(1030, 8929, 3438)
    5  = @(85, 155, 25)
    6  = @(25, 410, 139)
    7  = @(1, 0, 0)
    8  = @(341, 2441, 1025)
    10 = @(0, 0, 0)
    11 = @(0, 0, 0)
    12 = @(0, 6, 5)
    13 = @(0, 0, 0)
    14 = @(0, 0, 3)
    15 = @(0, 4, 0)
    16 = @(0, 1, 2)
    17 = @(3, 10, 1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}

# Row 9 only has Low/Moderate counts (no High value).
$ws.Cells.Item(9, 2).Value = 2
$ws.Cells.Item(9, 3).Value = 43

# Move the selection/active cell like the saved view did.
[void]$ws.Range("G18").Select()
